$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1. Add the new "ISSUE STATE" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$issueState = $wb.Worksheets.Add($null, $lastSheet)
$issueState.Name = "ISSUE STATE"

# 2. Build row 3 of the new sheet = exact copy of the row (MN) being removed from RESDT,
#    cell by cell so no extra blank cells get materialized.
$cols = @("A","B","C","D","E","F","I","J","K","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AU")
foreach ($col in $cols) {
    $src = $ws1.Range($col + "18")
    $dst = $issueState.Range($col + "3")
    $src.Copy($dst)
}

# 3. Build row 2 of the new sheet = new "CO" entry. Start from a full-row copy of the same
#    source row (so every column A:AU gets a cell, matching the target layout), then overwrite
#    the actual content/formatting that differs for the new state.
$issueState.Range("A2:AU2").ClearContents()
$ws1.Range("A18:AU18").Copy($issueState.Range("A2:AU2"))

$issueState.Range("A2").Value2 = "CO"
$issueState.Range("E2").Value2 = "White Test"
$issueState.Range("P2").Value2 = "123 Colorado CO usa"
$issueState.Range("AN2").Value2 = "52527"
$issueState.Range("AU2").Value2 = "Badger Test"

$issueState.Range("E2").Interior.Color = 14348258
$issueState.Range("F2").Interior.Color = 14348258

# 4. Remove the MN row from the RESDT sheet (shifts rows 19:42 up to 18:41)
$ws1.Rows.Item(18).Delete()

# 5. Make the new sheet the active tab with its own selection, matching the authored workbook.
$issueState.Select()
$issueState.Range("D11").Select()
